$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column divisors: F and N use 1000 (they hold values with 3 implied decimals),
# P uses 100 (2 implied decimals), and the rest (G,H,I,J,K,L,M,O) use 10 (1 implied decimal).
$divisors = @{
    "F" = 1000
    "G" = 10
    "H" = 10
    "I" = 10
    "J" = 10
    "K" = 10
    "L" = 10
    "M" = 10
    "N" = 1000
    "O" = 10
    "P" = 100
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in $divisors.Keys) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value2
        if ($current -ne $null) {
            $cell.Value2 = $current / $divisors[$col]
        }
    }
}
